$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before column A (shifts old A:D -> C:F)
$ws.Range("A1:B1").EntireColumn.Insert()

# Copy the header style (bold/centered/bordered "Pandas" style) from the
# existing header cell onto the two newly-inserted header cells.
$ws.Range("C1").Copy()
$ws.Range("A1:B1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# New header row 1 values
$ws.Range("A1").Value = "button_quantityIncrease_internalRoleButtonName"
$ws.Range("B1").Value = "button_quantityIncrease_nthChild"

# New data row 2 values (stored as text, matching the rest of the row).
# B2's value ("3") looks numeric, so assigning it through .Value/.Value2
# would auto-coerce it to a number; going through a text-returning formula
# and then converting the formula to its literal result keeps it text
# without requiring a dedicated (and residual) text number format/style.
$ws.Range("A2").Value = "+"
$ws.Range("B2").Formula = "=""3"""
$ws.Range("B2").Copy()
$ws.Range("B2").PasteSpecial(-4163)
$excel.CutCopyMode = $false

# Update the renamed path in what is now column C (was column A)
$ws.Range("C2").Value = "Data Files/AI-Generated/Common/fillShippingInfoAndCompleteOrder-test-data"

# Column widths (stored widths: 48, 34, 75, 9, 24, 21).
# COM ColumnWidth reads/writes with the usual Excel padding offset (~0.83
# less than the stored character width), matched against the existing
# columns (82 -> 81.17, 9 -> 8.17, 24 -> 23.17, 21 -> 20.17).
$ws.Columns.Item(1).ColumnWidth = 47.17
$ws.Columns.Item(2).ColumnWidth = 33.17
$ws.Columns.Item(3).ColumnWidth = 74.17
$ws.Columns.Item(4).ColumnWidth = 8.17
$ws.Columns.Item(5).ColumnWidth = 23.17
$ws.Columns.Item(6).ColumnWidth = 20.17
